# January-2021-Collection.xlsx — "Data updated till 7 Jan 11AM"
#
# Adds the day's collection figures (column L = 6-Jan-2021) for a batch of
# retailers, moves Vijay's breakdown comment from L66 to the now-populated
# L48, and fixes a retailer's display name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New collection amounts for 6-Jan (column L).
#    Rows whose other day-columns already carry the highlighted
#    ("s=11") look get that same look copied onto L; the rest just take
#    the sheet's plain default cell format for column L ("s=3"), which is
#    what a bare value write already produces.
# ---------------------------------------------------------------------

$highlightSource = $ws.Range("H3")   # an existing cell using the "highlighted" style

$highlightRows = @{
    7  = 7000
    11 = 2000
    19 = 2000
    25 = 3000
    48 = 6000
    73 = 3000
    77 = 2000
}

$plainRows = @{
    8  = 2000
    12 = 1000
    14 = 2500
    31 = 1000
    39 = 2000
    54 = 4000
    56 = 2000
    62 = 2000
    63 = 2000
    65 = 5000
    69 = 400
    70 = 3000
    71 = 2000
}

foreach ($row in $highlightRows.Keys) {
    $target = $ws.Range("L$row")
    $highlightSource.Copy()
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $target.Value = $highlightRows[$row]
}

foreach ($row in $plainRows.Keys) {
    $ws.Range("L$row").Value = $plainRows[$row]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Row 66 no longer carries any entry in column L — remove the cell
#    (and its comment, handled below) entirely rather than leaving it
#    blank.
# ---------------------------------------------------------------------

$ws.Range("L66").Comment.Delete()
$ws.Range("L66").Clear()

# ---------------------------------------------------------------------
# 3. Vijay's comment that used to sit on L66 now documents the new L48
#    figure instead.
# ---------------------------------------------------------------------

$ws.Range("L48").AddComment("Vijay:`n1700-Cash`n2300-Digital`n2000-Digital") | Out-Null

# ---------------------------------------------------------------------
# 4. Retailer display-name correction.
# ---------------------------------------------------------------------

$ws.Range("B73").Value = "INDRJEET KUMAR"

# ---------------------------------------------------------------------
# 5. Bring the view roughly in line with where work left off: frozen
#    panes scrolled down near the newest rows, with L77 selected.
# ---------------------------------------------------------------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("G3").Select()
$win.FreezePanes = $true
$ws.Range("L77").Select()
